$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 11
$ws.Range("E2").Value = 10.65

# Row 3
$ws.Range("D3").Value = 10.5

# Row 4
$ws.Range("B4").Value = 9
$ws.Range("C4").Value = 9.5
$ws.Range("E4").Value = 10.67
$ws.Range("F4").Value = 10.11

# Row 5
$ws.Range("B5").Value = 9.35
$ws.Range("D5").Value = 9.33
$ws.Range("F5").Value = 10.35
$ws.Range("H5").Value = 8.83

# Row 6
$ws.Range("D6").Value = 9.890000000000001
$ws.Range("E6").Value = 9.65
$ws.Range("G6").Value = 10.27
$ws.Range("H6").Value = 10.74

# Row 7
$ws.Range("F7").Value = 9.73
$ws.Range("J7").Value = 9.380000000000001

# Row 8
$ws.Range("E8").Value = 11.17
$ws.Range("F8").Value = 9.26
$ws.Range("J8").Value = 11.71

# Row 10
$ws.Range("G10").Value = 10.62
$ws.Range("H10").Value = 8.289999999999999
